# Apply trade #115 close update across the workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1198.85   # Current Capital
$summary.Range("B4").Value = -1.16     # Total P&L $
$summary.Range("B5").Value = -0.2      # Total P&L %
$summary.Range("B6").Value = 115       # Total Trades
$summary.Range("B7").Value = 41        # Winning Trades
$summary.Range("B9").Value = 35.65     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 98.84999999999999  # Capital
$status.Range("D4").Value = 115                # Trades
$status.Range("E4").Value = -1.16              # P&L $
$status.Range("F4").Value = -1.15              # P&L %
$status.Range("G4").Value = 35.65              # Win Rate %

# ---------------------------------------------------------------------
# Helper to append the new trade row (#115 -> row 116) to a trades sheet
# ---------------------------------------------------------------------
function Add-TradeRow115 {
    param($ws)

    $ws.Cells.Item(116, 1).Value = 115

    # Column B holds the date as plain text ("2026-02-17"), matching every
    # other row in the sheet. Mark the cell as Text first so Excel doesn't
    # auto-convert the literal into a date serial number.
    $ws.Cells.Item(116, 2).NumberFormat = "@"
    $ws.Cells.Item(116, 2).Value = "2026-02-17"

    $ws.Cells.Item(116, 3).Value = "16:03:10"
    $ws.Cells.Item(116, 4).Value = "MarketMaking"
    $ws.Cells.Item(116, 5).Value = "UP"
    $ws.Cells.Item(116, 6).Value = 0.84
    $ws.Cells.Item(116, 7).Value = 0.85
    $ws.Cells.Item(116, 8).Value = "CLOSED"
    $ws.Cells.Item(116, 9).Value = 1.1905
    $ws.Cells.Item(116, 10).Value = 0.01
    $ws.Cells.Item(116, 11).Value = 98.84999999999999
    $ws.Cells.Item(116, 12).Value = 0
    $ws.Cells.Item(116, 13).Value = 0
    $ws.Cells.Item(116, 14).Value = 0.6
    $ws.Cells.Item(116, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(116, 16).Value = "early_exit"
    $ws.Cells.Item(116, 17).Value = 0.14
}

# ---------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow115 $allTrades

# ---------------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow115 $marketMaking
